$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.867.24"
$ws.Range("E2").Value = "  +2.40%  "
$ws.Range("D3").Value = "2.587.81"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'520.36"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("D6").Value = "'139.74"
$ws.Range("E6").Value = "  -2.73%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +0.68%  "
$ws.Range("D9").Value = "2.600.36"
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("D10").Value = "'6.54"
$ws.Range("E10").Value = "  -0.55%  "
$ws.Range("E11").Value = "  +0.57%  "
$ws.Range("E12").Value = "  +1.67%  "
$ws.Range("E13").Value = "  +2.92%  "
$ws.Range("D14").Value = "3.044.98"
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("D15").Value = "58.836.47"
$ws.Range("E15").Value = "  +2.36%  "
$ws.Range("E16").Value = "  +1.67%  "
$ws.Range("D17").Value = "2.603.63"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").Value = "'338.95"
$ws.Range("E19").Value = "  +1.17%  "
$ws.Range("E20").Value = "  +0.45%  "
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("E22").Value = "  +3.56%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("E24").Value = "  +2.43%  "
$ws.Range("E25").Value = "  +0.98%  "
$ws.Range("E26").Value = "  +1.03%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("E28").Value = "  +1.64%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").Value = "0.0₃0723"
$ws.Range("E30").Value = "  -3.51%  "
$ws.Range("E31").Value = "  -4.89%  "
$ws.Range("E32").Value = "  -1.13%  "
$ws.Range("E33").Value = "  +1.06%  "
$ws.Range("D34").Value = "'149.00"
$ws.Range("E34").Value = "  +0.39%  "
$ws.Range("E35").Value = "  -1.06%  "
$ws.Range("D36").Value = "'1.12"
$ws.Range("E36").Value = "  -1.18%  "
$ws.Range("D37").Value = "'36.32"
$ws.Range("E37").Value = "  +1.34%  "
$ws.Range("E38").Value = "  +2.59%  "
$ws.Range("E39").Value = "  -0.63%  "
$ws.Range("D40").Value = "'0.819"
$ws.Range("E40").Value = "  -2.21%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").Value = "'0.998"
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").Value = "'274.10"
$ws.Range("E43").Value = "  +2.48%  "
$ws.Range("D44").Value = "'10.76"
$ws.Range("E44").Value = "  +1.02%  "
$ws.Range("D45").Value = "'0.591"
$ws.Range("E45").Value = "  +0.41%  "
$ws.Range("E46").Value = "  -0.16%  "
$ws.Range("D47").Value = "'0.0521"
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("D48").Value = "'18.57"
$ws.Range("E48").Value = "  -1.19%  "
$ws.Range("D49").Value = "1.989.51"
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("E50").Value = "  +0.35%  "
$ws.Range("E51").Value = "  -1.31%  "
